# Handle exist category in db: append two new category values ("test", "test2")
# to the categoryval sheet, growing the used range from A1:A5 to A1:A7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "test"
$ws.Range("A7").Value = "test2"

# Move/update the active selection to the last appended cell (A7),
# matching Excel's natural behavior after data entry.
$ws.Range("A7").Select()
